# Michigan_Basis.xlsx — "Updated policies and graphs"
# Fill in newly-tracked start/end policy dates for Casinos (row 28) and
# Bars for indoor service / Food and Drink (row 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 (Casinos): add start_date (C28)
$ws.Range("C28").Value = 43906
$ws.Range("C28").NumberFormat = "YYYY-MM-DD"

# Row 34 (Bars for indoor service): add start_date (C34) and end_dates (D34)
$ws.Range("C34").Value = 43906
$ws.Range("C34").NumberFormat = "YYYY-MM-DD"
$ws.Range("D34").Value = 43990
$ws.Range("D34").NumberFormat = "YYYY-MM-DD"

# Reflect the author's final scroll position / selection when they saved.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("D35").Select()
